$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing "Data" sheet (transactional data) to
#    "Transactional Data" BEFORE inserting the new sheet, so that the name
#    "Data" is free for the newly inserted sheet.
# ---------------------------------------------------------------------------
$acct    = $wb.Worksheets.Item("Account")
$oldData = $wb.Worksheets.Item("Data")
$oldData.Name = "Transactional Data"

# ---------------------------------------------------------------------------
# 2. Insert a brand-new worksheet right after "Account" (i.e. before
#    "Transactional Data") and name it "Data". This becomes the new
#    "Account dimension metadata" sheet.
# ---------------------------------------------------------------------------
$newData = $wb.Worksheets.Add($null, $acct)
$newData.Name = "Data"

# ---------------------------------------------------------------------------
# 3. Populate the header row (row 1) with bold formatting.
# ---------------------------------------------------------------------------
$headers = @(
    "Member ID",
    "Description",
    "Hierarchy",
    "Formula",
    "Account Type",
    "Calculated On",
    "Aggregation Type",
    "Excepion Aggregation Type",
    "Exception Aggregation Dimension",
    "Required Dimensions",
    "Scale",
    "Decimal Places",
    "Units & Currencies",
    "Thresholds",
    "Hide"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newData.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# ---------------------------------------------------------------------------
# 4. Populate the single data row (row 2).
# ---------------------------------------------------------------------------
$newData.Range("A2").Value = "DATA"
$newData.Range("B2").Value = "Data"
$newData.Range("C2").Value = "<root>"
$newData.Range("E2").Value = "NFIN"
$newData.Range("F2").Value = "SUM"
$newData.Range("O2").Value = $false

# ---------------------------------------------------------------------------
# 5. Size the columns to fit their content (matches the "bestFit" columns
#    seen in the target workbook).
# ---------------------------------------------------------------------------
$newData.Range("A1:O2").EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 6. Selection / active cell bookkeeping.
# ---------------------------------------------------------------------------
$newData.Activate()
$newData.Range("H8").Select() | Out-Null
